# Update gh-pages to output generated at 456a3b4
# Applies updated "F" column (participants/views count) values across the
# 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 860
$ws1.Range("F3").Value  = 13790
$ws1.Range("F4").Value  = 13577
$ws1.Range("F12").Value = 762
$ws1.Range("F14").Value = 98
$ws1.Range("F15").Value = 92
$ws1.Range("F20").Value = 433
$ws1.Range("F21").Value = 399
$ws1.Range("F22").Value = 323
$ws1.Range("F23").Value = 263
$ws1.Range("F24").Value = 834
$ws1.Range("F25").Value = 86

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 48
$ws2.Range("F7").Value = 1502

# Sheet 4: 全部类型 (aggregated view of all the above rows, offset by 1 row)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 860
$ws4.Range("F4").Value  = 13790
$ws4.Range("F5").Value  = 13577
$ws4.Range("F13").Value = 762
$ws4.Range("F14").Value = 48
$ws4.Range("F17").Value = 98
$ws4.Range("F18").Value = 92
$ws4.Range("F27").Value = 433
$ws4.Range("F28").Value = 399
$ws4.Range("F29").Value = 323
$ws4.Range("F30").Value = 263
$ws4.Range("F31").Value = 834
$ws4.Range("F33").Value = 1502
$ws4.Range("F37").Value = 86

$wb.Save()
